$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "TODO" test row (row 43), matching the surrounding
# Train / Branching / Train Dialog rows already in the grid.
$ws.Range("A43").Value = "Train"
$ws.Range("B43").Value = "Branching"
$ws.Range("C43").Value = "Train Dialog"
$ws.Range("D43").Value = "TODO: Branching + Edit how do they mix? Edit first does not allow branching."

# Highlight D43 as a TODO note: wrap text, red font, yellow fill.
$d43 = $ws.Range("D43")
$d43.WrapText = $true
$d43.Font.Color = 255
$d43.Interior.Color = 65535

# Move the view/selection down to the newly added row.
$win = $ws.Application.ActiveWindow
$win.ScrollRow = 32
$ws.Range("D43").Select() | Out-Null
